# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the first data row (the d2f80547... file) on both the
# "zh-cn" and "de-de" report sheets, to reflect a freshly generated
# handback report.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 09:39:16"
$wsZh.Range("H2").Value = "2016-03-23 09:39:56"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 09:39:23"
$wsDe.Range("H2").Value = "2016-03-23 09:40:12"
